$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -0.00000003973643103449831
# Row 3
$ws.Range("C3").Value = 11.03038649796165
$ws.Range("D3").Value = 9.159267103524588
$ws.Range("E3").Value = 0.1597303064348723
$ws.Range("F3").Value = 2.713836019108593
$ws.Range("G3").Value = 1.260165004315388
$ws.Range("H3").Value = 1.357485128991759
$ws.Range("I3").Value = 1.860954346686789
$ws.Range("J3").Value = 1.184153455403833
$ws.Range("K3").Value = 6.20458584844114
$ws.Range("L3").Value = 0.01252840826584656
$ws.Range("M3").Value = 0.005933893446967819
# Row 4
$ws.Range("C4").Value = 11.03038649796165
$ws.Range("D4").Value = 9.159267103524588
$ws.Range("E4").Value = -0.05634054514634194
$ws.Range("F4").Value = 3.039639152562313
$ws.Range("G4").Value = 1.220843741232643
$ws.Range("H4").Value = 1.495549088943689
$ws.Range("I4").Value = 1.902072993280568
$ws.Range("J4").Value = 1.109834445167124
$ws.Range("K4").Value = 6.20458584844114
$ws.Range("L4").Value = 0.01252840826584656
$ws.Range("M4").Value = 0.005933893446967819
# Row 5
$ws.Range("C5").Value = 11.03038649796165
$ws.Range("D5").Value = 9.159267103524588
$ws.Range("E5").Value = -0.05634054514634194
$ws.Range("F5").Value = 3.039639152562313
$ws.Range("G5").Value = 1.220843741232643
$ws.Range("H5").Value = 1.495549088943689
$ws.Range("I5").Value = 1.902072993280568
$ws.Range("J5").Value = 1.109834445167124
$ws.Range("K5").Value = 6.20458584844114
$ws.Range("L5").Value = 0.01252840826584656
$ws.Range("M5").Value = 0.005933893446967819
# Row 6
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 11.03038649796165
$ws.Range("D6").Value = 9.159267103524588
$ws.Range("E6").Value = -0.03129893233941866
$ws.Range("F6").Value = 3.839575757841913
$ws.Range("G6").Value = -0.01994980750244963
$ws.Range("H6").Value = 0.0000000008003052725756468
$ws.Range("I6").Value = 0.00000000000000000000247658767146454
$ws.Range("J6").Value = 0.000000000000000000000114057367153024
$ws.Range("K6").Value = 3.424903828327841
$ws.Range("L6").Value = 0.000000000000000000000000000133317492982351
$ws.Range("M6").Value = 0.0000000000000000000000000001336035820230856
# Row 7
$ws.Range("C7").Value = 11.03038649796165
$ws.Range("D7").Value = 9.159267103524588
$ws.Range("E7").Value = 0.008651569081671613
$ws.Range("F7").Value = 3.580009255957597
$ws.Range("G7").Value = 0.1127101358838678
$ws.Range("H7").Value = 0.0000000000000000000000000004406086709259265
$ws.Range("I7").Value = 0.000000000000000000000000000000000000000000000000000001363339214312259
$ws.Range("J7").Value = 0.0000000000000000000000000000000000000000000000000000003322689427154167
$ws.Range("K7").Value = 3.424834880236386
$ws.Range("L7").Value = 0.000000000000000000000000001119262147691079
$ws.Range("M7").Value = 0.000000000000000000000000001092007793405155
# Row 8
$ws.Range("C8").Value = 11.03038649796165
$ws.Range("D8").Value = 9.159267103524588
$ws.Range("E8").Value = 0.1202246781991047
$ws.Range("F8").Value = 2.900803018595902
$ws.Range("G8").Value = 0.8386510598470626
$ws.Range("H8").Value = 0.5118128579440012
$ws.Range("I8").Value = 0.001108584282262656
$ws.Range("J8").Value = 0.0002461425776835096
$ws.Range("K8").Value = 5.068387254402111
$ws.Range("L8").Value = 0.000000000000000000000000000133317492982351
$ws.Range("M8").Value = 0.0000000000000000000000000001336035820230856
# Row 9
$ws.Range("C9").Value = 11.03038649796165
$ws.Range("D9").Value = 9.159267103524588
$ws.Range("E9").Value = 0.005113958203569938
$ws.Range("F9").Value = 11.84572237979037
$ws.Range("G9").Value = -0.6580550343440751
$ws.Range("H9").Value = 2.391212601116135
$ws.Range("I9").Value = 12.15337514323337
$ws.Range("J9").Value = 0.6832001467173409
$ws.Range("K9").Value = 3.424834880236386
$ws.Range("L9").Value = 0.000000000000000000000000001109335647967048
$ws.Range("M9").Value = 0.000000000000000000000000001082268069760972
# Row 10
$ws.Range("C10").Value = 11.03038649796165
$ws.Range("D10").Value = 9.159267103524588
$ws.Range("E10").Value = 0.006124853759927205
$ws.Range("F10").Value = 9.877710933824162
$ws.Range("G10").Value = -1.184622164985516
$ws.Range("H10").Value = 2.217541146533207
$ws.Range("I10").Value = 8.79065527282544
$ws.Range("J10").Value = -0.6359761225589421
$ws.Range("K10").Value = 3.424834880236386
$ws.Range("L10").Value = 0.000000000000000000000000001109335647967048
$ws.Range("M10").Value = 0.000000000000000000000000001082268069760972
# Row 11
$ws.Range("C11").Value = 11.03038649796165
$ws.Range("D11").Value = 9.159267103524588
$ws.Range("E11").Value = 0.006124853759927205
$ws.Range("F11").Value = 9.877710933824162
$ws.Range("G11").Value = -1.184622164985516
$ws.Range("H11").Value = 2.217541146533207
$ws.Range("I11").Value = 8.79065527282544
$ws.Range("J11").Value = -0.6359761225589421
$ws.Range("K11").Value = 3.424834880236386
$ws.Range("L11").Value = 0.000000000000000000000000001109335647967048
$ws.Range("M11").Value = 0.000000000000000000000000001082268069760972
# Row 15
$ws.Range("E15").Value = 0.0000304982991176699
$ws.Range("F15").Value = 0.000001322257755762889
$ws.Range("G15").Value = 0.0000001687078097113069
$ws.Range("H15").Value = 0.000002853234980418634
$ws.Range("I15").Value = 0.00000000001155619596986963
$ws.Range("J15").Value = 0.00000000001049031921478155
$ws.Range("K15").Value = 0.000005735841689583602
$ws.Range("L15").Value = 0.00000000002524842922867389
$ws.Range("M15").Value = 0.00000000002519518469257158
# Row 16
$ws.Range("E16").Value = -0.0002144720849258486
$ws.Range("F16").Value = 0.000008336368421167271
$ws.Range("G16").Value = 0.000001316684930640294
$ws.Range("H16").Value = 0.00002168190164791682
$ws.Range("I16").Value = 0.000000001922994564719176
$ws.Range("J16").Value = 0.000000001758373793877078
$ws.Range("K16").Value = 0.00004033874748846599
$ws.Range("L16").Value = 0.000000001945665935367659
$ws.Range("M16").Value = 0.000000001943304139872167
# Row 17
$ws.Range("E17").Value = -0.0002144720849258486
$ws.Range("F17").Value = 0.000008336368421167271
$ws.Range("G17").Value = 0.000001316684930640294
$ws.Range("H17").Value = 0.00002168190164791682
$ws.Range("I17").Value = 0.000000001922994564719176
$ws.Range("J17").Value = 0.000000001758373793877078
$ws.Range("K17").Value = 0.00004033874748846599
$ws.Range("L17").Value = 0.000000001945665935367659
$ws.Range("M17").Value = 0.000000001943304139872167
# Row 18
$ws.Range("E18").Value = 0.007048338896217496
$ws.Range("F18").Value = 0.0001527339526123465
$ws.Range("G18").Value = 0.0001381315171326805
$ws.Range("H18").Value = 0.000000184446094366096
$ws.Range("I18").Value = 0.00000000000001384681615131152
$ws.Range("J18").Value = 0.00000000000001383670946399082
$ws.Range("K18").Value = 0.0003486512894174307
$ws.Range("L18").Value = 0.00000004697690863719423
$ws.Range("M18").Value = 0.00000004701483017177481
# Row 19
$ws.Range("E19").Value = 0.007075978366416692
$ws.Range("F19").Value = 0.00002263211062356349
$ws.Range("G19").Value = 0.00002133843775716282
$ws.Range("H19").Value = 0.00000004955107630319862
$ws.Range("I19").Value = 0.00000000000000018722767588752
$ws.Range("J19").Value = 0.0000000000000001856039591028094
$ws.Range("K19").Value = 0.00009452878061942556
$ws.Range("L19").Value = 0.0000000001357640163417291
$ws.Range("M19").Value = 0.0000000001348201814836058
# Row 20
$ws.Range("E20").Value = 0.007075978366416692
$ws.Range("F20").Value = 0.00002263211062356349
$ws.Range("G20").Value = 0.00002133843775716282
$ws.Range("H20").Value = 0.00000004955107630319862
$ws.Range("I20").Value = 0.00000000000000018722767588752
$ws.Range("J20").Value = 0.0000000000000001856039591028094
$ws.Range("K20").Value = 0.00009452878061942556
$ws.Range("L20").Value = 0.0000000001357640163417291
$ws.Range("M20").Value = 0.0000000001348201814836058
# Row 21
$ws.Range("E21").Value = 0.000001959804056648239
$ws.Range("F21").Value = 0.000004298119671357475
$ws.Range("G21").Value = -0.000001383309104355598
$ws.Range("H21").Value = 0.0000003769088239556561
$ws.Range("I21").Value = 0.0000000000004000730894164373
$ws.Range("J21").Value = 0.00000000000004057943042991756
$ws.Range("K21").Value = 0.0000009861348368513163
$ws.Range("L21").Value = 0.0000000000005271151386122459
$ws.Range("M21").Value = 0.0000000000005038616579524655
# Row 22
$ws.Range("E22").Value = 0.0001860131251966351
$ws.Range("F22").Value = 0.000006887687052059391
$ws.Range("G22").Value = -0.0000009563522462812285
$ws.Range("H22").Value = 0.000002225745641289322
$ws.Range("I22").Value = 0.00000000001676813647122356
$ws.Range("J22").Value = 0.000000000003515764797321014
$ws.Range("K22").Value = 0.000002015933335195811
$ws.Range("L22").Value = 0.000000000004060592043585097
$ws.Range("M22").Value = 0.000000000003945789335512037
# Row 23
$ws.Range("E23").Value = 0.0001860131251966351
$ws.Range("F23").Value = 0.000006887687052059391
$ws.Range("G23").Value = -0.0000009563522462812285
$ws.Range("H23").Value = 0.000002225745641289322
$ws.Range("I23").Value = 0.00000000001676813647122356
$ws.Range("J23").Value = 0.000000000003515764797321014
$ws.Range("K23").Value = 0.000002015933335195811
$ws.Range("L23").Value = 0.000000000004060592043585097
$ws.Range("M23").Value = 0.000000000003945789335512037
# Row 24
$ws.Range("E24").Value = 0.01143641912412727
$ws.Range("F24").Value = 0.00008642161448091815
$ws.Range("G24").Value = 0.00007756925547959815
$ws.Range("H24").Value = 0.00008736292356318004
$ws.Range("I24").Value = 0.0000000021580709279565
$ws.Range("J24").Value = 0.000000002105242540914173
$ws.Range("K24").Value = 0.0003235792745324306
$ws.Range("L24").Value = 0.00000002273753585088633
$ws.Range("M24").Value = 0.00000002219779971211617
# Row 25
$ws.Range("E25").Value = 0.01664945516026482
$ws.Range("F25").Value = 0.0001693705925941238
$ws.Range("G25").Value = 0.000149537058998303
$ws.Range("H25").Value = 0.00000141373907746232
$ws.Range("I25").Value = 0.000000000000730047368109126
$ws.Range("J25").Value = 0.0000000000007126998999441604
$ws.Range("K25").Value = 0.0002391891593518394
$ws.Range("L25").Value = 0.00000001847221489921052
$ws.Range("M25").Value = 0.00000001804198138130919
# Row 26
$ws.Range("E26").Value = 0.01280264361116186
$ws.Range("F26").Value = 0.00004630123976433341
$ws.Range("G26").Value = 0.00004196067875640344
$ws.Range("H26").Value = 0.00002231439456564559
$ws.Range("I26").Value = 0.000000000164801397655328
$ws.Range("J26").Value = 0.0000000001609346014572481
$ws.Range("K26").Value = 0.00018203289845917
$ws.Range("L26").Value = 0.000000008168899662019224
